$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.515.50"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "2.109.56"
$ws.Range("E3").Value = "  -0.31%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'335.22"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").Value = "'0.5239"
$ws.Range("E7").Value = "  -1.54%  "

$ws.Range("D8").Value = "'0.4546"
$ws.Range("E8").Value = "  +3.69%  "

$ws.Range("D9").Value = "'53.33"
$ws.Range("E9").Value = "  +15.74%  "

$ws.Range("D10").Value = "'0.09010"
$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("D11").Value = "'1.162"
$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("D12").Value = "'24.53"
$ws.Range("E12").Value = "  -2.15%  "

$ws.Range("D13").Value = "2.110.58"
$ws.Range("E13").Value = "  -0.20%  "

$ws.Range("D14").Value = "'6.788"
$ws.Range("E14").Value = "  +0.14%  "

$ws.Range("D15").Value = "'7.867"
$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").Value = "'0.00001127"
$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("D19").Value = "'0.06625"
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").Value = "'19.35"
$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "'6.309"
$ws.Range("E22").Value = "  -0.62%  "

$ws.Range("D23").Value = "30.565.25"
$ws.Range("E23").Value = "  -1.06%  "

$ws.Range("D24").Value = "'12.39"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").Value = "'2.350"
$ws.Range("E25").Value = "  +3.46%  "

$ws.Range("D26").Value = "2.351.74"
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("D27").Value = "'22.43"
$ws.Range("E27").Value = "  -1.72%  "

$ws.Range("D28").Value = "'2.575"
$ws.Range("E28").Value = "  -0.49%  "

$ws.Range("D29").Value = "'163.41"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "'133.11"
$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").Value = "'1.203"
$ws.Range("E31").Value = "  +1.12%  "

$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("D33").Value = "'1.669"
$ws.Range("E33").Value = "  +6.74%  "

$ws.Range("D34").Value = "'6.170"
$ws.Range("E34").Value = "  -1.37%  "

$ws.Range("D35").Value = "'3.952"
$ws.Range("E35").Value = "  -1.56%  "

$ws.Range("D36").Value = "'10.58"
$ws.Range("E36").Value = "  +11.39%  "

$ws.Range("D37").Value = "'0.02583"
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("D38").Value = "'0.06810"
$ws.Range("E38").Value = "  +0.57%  "

$ws.Range("D39").Value = "'5.544"
$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").Value = "'12.76"
$ws.Range("E40").Value = "  -1.20%  "

$ws.Range("D41").Value = "'0.2287"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").Value = "'0.6929"
$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").Value = "'2.395"
$ws.Range("E44").Value = "  +7.17%  "

$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.20%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6412"
$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'14.05"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D48").Value = "'3.668"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").Value = "'1.249"
$ws.Range("E49").Value = "  -2.13%  "

$ws.Range("E50").Value = "  +5.30%  "

$ws.Range("D51").Value = "'83.46"
$ws.Range("E51").Value = "  +0.54%  "
